$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "Sum total: x1:1.000|x2:1.000|x3:2.000|x4:1.000|x5:1.000|x6:1.000|x7:0.563|x8:1.000|x9:1.000|x10:1.000|x11:1.000|x12:0.437|x13:1.000|x14:1.000|x15:1.000|x16:1.000|x17:0.000"
$ws.Range("I2").Value = "0.941, 0.380"
$ws.Range("J2").Value = "88.172, 89.247"
$ws.Range("K2").Value = "99.810, 99.810"
